$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '46.342.65'
$ws.Range("E2").Value = '  -0.82%  '

# Row 3
$ws.Range("D3").Value = '2.468.47'
$ws.Range("E3").Value = '  +8.47%  '

# Row 4
$ws.Range("E4").Value = '  +0.05%  '

# Row 5
$ws.Range("D5").Value = "'298.04"
$ws.Range("E5").Value = '  -0.84%  '

# Row 6
$ws.Range("D6").Value = "'97.00"
$ws.Range("E6").Value = '  -2.73%  '

# Row 7
$ws.Range("D7").Value = "'0.579"
$ws.Range("E7").Value = '  +0.71%  '

# Row 8
$ws.Range("E8").Value = '  +0.01%  '

# Row 9
$ws.Range("E9").Value = '  +1.49%  '

# Row 10
$ws.Range("D10").Value = "'35.53"
$ws.Range("E10").Value = '  +1.32%  '

# Row 11
$ws.Range("D11").Value = "'0.0790"
$ws.Range("E11").Value = '  -1.93%  '

# Row 12
$ws.Range("E12").Value = '  +2.02%  '

# Row 13
$ws.Range("E13").Value = '  +2.04%  '

# Row 14
$ws.Range("D14").Value = '2.836.99'
$ws.Range("E14").Value = '  +8.18%  '

# Row 15
$ws.Range("D15").Value = '2.454.14'
$ws.Range("E15").Value = '  +8.03%  '

# Row 16
$ws.Range("D16").Value = "'0.862"
$ws.Range("E16").Value = '  +7.69%  '

# Row 17
$ws.Range("D17").Value = "'14.21"
$ws.Range("E17").Value = '  +4.23%  '

# Row 18
$ws.Range("D18").Value = '46.382.92'
$ws.Range("E18").Value = '  -0.58%  '

# Row 19
$ws.Range("D19").Value = "'12.89"
$ws.Range("E19").Value = '  +3.08%  '

# Row 20
$ws.Range("D20").Value = '0.0₃0952'
$ws.Range("E20").Value = '  -4.42%  '

# Row 21
$ws.Range("D21").Value = "'6.31"
$ws.Range("E21").Value = '  +8.00%  '

# Row 22
$ws.Range("D22").Value = "'67.85"
$ws.Range("E22").Value = '  +2.81%  '

# Row 23
$ws.Range("D23").Value = "'246.85"
$ws.Range("E23").Value = '  -0.24%  '

# Row 24
$ws.Range("D24").Value = "'2.83"
$ws.Range("E24").Value = '  +1.43%  '

# Row 25
$ws.Range("D25").Value = "'1.99"
$ws.Range("E25").Value = '  +6.64%  '

# Row 26
$ws.Range("E26").Value = '  -0.07%  '

# Row 27
$ws.Range("D27").Value = "'40.38"
$ws.Range("E27").Value = '  -2.10%  '

# Row 28
$ws.Range("E28").Value = '  -1.20%  '

# Row 29
$ws.Range("D29").Value = "'9.89"
$ws.Range("E29").Value = '  +2.98%  '

# Row 30
$ws.Range("D30").Value = "'3.94"
$ws.Range("E30").Value = '  +17.16%  '

# Row 31
$ws.Range("D31").Value = "'21.64"
$ws.Range("E31").Value = '  +7.47%  '

# Row 32
$ws.Range("D32").Value = "'5.63"
$ws.Range("E32").Value = '  +5.10%  '

# Row 33
$ws.Range("D33").Value = "'2.77"
$ws.Range("E33").Value = '  -1.30%  '

# Row 34
$ws.Range("D34").Value = "'148.81"
$ws.Range("E34").Value = '  +1.49%  '

# Row 35
$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").Value = "'2.07"
$ws.Range("E35").Value = '  +23.03%  '

# Row 36
$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D36").Value = "'0.0780"
$ws.Range("E36").Value = '  +1.42%  '

# Row 37
$ws.Range("E37").Value = '  +1.95%  '

# Row 38
$ws.Range("E38").Value = '  +0.49%  '

# Row 39
$ws.Range("D39").Value = "'15.42"
$ws.Range("E39").Value = '  -1.03%  '

# Row 40
$ws.Range("D40").Value = "'4.00"
$ws.Range("E40").Value = '  +3.93%  '

# Row 41
$ws.Range("E41").Value = '  +2.30%  '

# Row 42
$ws.Range("E42").Value = '  +7.34%  '

# Row 43
$ws.Range("D43").Value = '1.997.73'
$ws.Range("E43").Value = '  +11.95%  '

# Row 44
$ws.Range("E44").Value = '  -0.02%  '

# Row 45
$ws.Range("D45").Value = "'92.53"
$ws.Range("E45").Value = '  -1.53%  '

# Row 46
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").Value = "'16.69"
$ws.Range("E46").Value = '  +34.61%  '

# Row 47
$ws.Range("B47").Value = 'Stacks'
$ws.Range("C47").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D47").Value = "'1.82"
$ws.Range("E47").Value = '  -3.53%  '

# Row 48
$ws.Range("D48").Value = "'8.67"
$ws.Range("E48").Value = '  +9.87%  '

# Row 49
$ws.Range("D49").Value = "'102.65"
$ws.Range("E49").Value = '  +8.18%  '

# Row 50
$ws.Range("D50").Value = '2.700.67'
$ws.Range("E50").Value = '  +8.08%  '

# Row 51
$ws.Range("E51").Value = '  +1.76%  '
